$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "step",
  "favor",
  "feel",
  "rather",
  "polici",
  "thought",
  "slight",
  "period",
  "rate",
  "governor",
  "spade",
  "set",
  "percent",
  "last",
  "term",
  "mayb",
  "higher",
  "fed",
  "mean",
  "growth",
  "use",
  "posit",
  "someth",
  "open",
  "degre",
  "can",
  "trade",
  "number",
  "issu",
  "maintain",
  "look",
  "discuss",
  "accept",
  "toward",
  "thing",
  "don",
  "half",
  "ii",
  "short",
  "work",
  "end",
  "prefer",
  "sensit",
  "remain",
  "statement",
  "call",
  "rang",
  "current",
  "intern",
  "desk",
  "monetari",
  "committe",
  "sentenc",
  "borrow",
  "reason",
  "peter",
  "guess",
  "need",
  "right",
  "iii",
  "concern",
  "exist",
  "clear",
  "know",
  "risk",
  "agre",
  "much",
  "kind",
  "fund",
  "appropri",
  "indic",
  "object",
  "longer",
  "realli",
  "give",
  "vote",
  "situat",
  "reserv",
  "eas",
  "seem",
  "side",
  "continu",
  "impli",
  "read",
  "understand",
  "anyth",
  "year",
  "certain",
  "mr",
  "week",
  "million",
  "associ",
  "alreadi",
  "will",
  "day",
  "discount",
  "suggest",
  "reduct",
  "problem",
  "presid",
  "talk",
  "market",
  "recent",
  "yes",
  "sure",
  "tighten",
  "tri",
  "show",
  "economi",
  "consist",
  "want",
  "move",
  "chairman",
  "time",
  "point",
  "condit",
  "probabl",
  "relationship",
  "question",
  "lower",
  "actual",
  "target",
  "normal",
  "financi",
  "pressur",
  "run",
  "given",
  "done",
  "littl",
  "altern",
  "convey",
  "around",
  "come",
  "languag",
  "view",
  "might",
  "uncertainti",
  "think",
  "direct",
  "chang",
  "meet",
  "let",
  "fact",
  "relat",
  "oper",
  "flexibl",
  "place",
  "expect",
  "differ",
  "feder"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

"done"